# Fixes a handful of typos / wording tweaks scattered across the deck, plus
# shrinks the GitHub link on the last slide.
#
# Each text edit below locates the run by its *exact, full* original text
# (via TextRange.Find) and overwrites just that run's text - this keeps
# every other run/paragraph in the shape completely untouched, matching
# the minimal nature of the source diff.

$p = $ppt.ActivePresentation

# --- Slide 2 : "Коммерческие реализации" column ---------------------------
$s2 = $p.Slides.Item(2)

$run = $s2.Shapes.Item(2).TextFrame.TextRange.Find("Коммерчиские реализации:")
$run.Text = "Коммерческие реализации:"

# --- Slide 2 : "Теоретические" column --------------------------------------
$run = $s2.Shapes.Item(3).TextFrame.TextRange.Find("Теоритические:")
$run.Text = "Теоретические:"

$run = $s2.Shapes.Item(3).TextFrame.TextRange.Find("интерпретируемый ")
$run.Text = "интерпретируемой "

# --- Slide 2 : "Open source" column -----------------------------------------
$run = $s2.Shapes.Item(4).TextFrame.TextRange.Find("Качество ниже чем у коммерческих пректов")
$run.Text = "Качество ниже чем у коммерческих проектов"

# --- Slide 6 : узловые точки ------------------------------------------------
$s6 = $p.Slides.Item(6)
$run = $s6.Shapes.Item(3).TextFrame.TextRange.Find("Например количество точек можно взять как квадратный корень из суммы разрешений изображения. ")
$run.Text = "Например количество точек можно взять как квадратный корень из суммы разрешения изображения. "

# --- Slide 7 : области не смежные с краями ----------------------------------
$s7 = $p.Slides.Item(7)
$run = $s7.Shapes.Item(1).TextFrame.TextRange.Find("  Берем только области, не смежные с краями изображения, и считаем разность. Размер матрицы снова уменьшается. ")
$run.Text = "  Берем только области, не смежные с краями изображения, и считаем их разность между собой. Размер матрицы снова уменьшается. "

# --- Slide 9 : shrink the GitHub link text ----------------------------------
$s9 = $p.Slides.Item(9)
$run = $s9.Shapes.Item(7).TextFrame.TextRange.Find("https://github.com/nikit34/Etsy_analys_MY_ALGORITHM_DETECT_DIFFERENCE_IMG")
$run.Font.Size = 14
